$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.247.70"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.939.63"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "2.935.67"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "65.242.38"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "3.429.41"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "2.939.72"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.44%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.973"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  -9.63%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -7.75%  "
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "385.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "2.704.98"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +4.67%  "
$ws.Range("E51").Value = "  -0.51%  "
